$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: "highlighs" -> "highlights" (M1), also picks up a plain black Arial,
# non-wrapping look (font formatting reset during the edit)
$ws.Range("M1").Style = $ws.Range("A1").Style
$ws.Range("M1").Font.Color = 0
$ws.Range("M1").Value = "highlights"

# Product 1 (row 2): rename color variant
$ws.Range("C2").Value = "camisa laranja"

# Product 2 (row 3): new name, price, desc_price, route and availability
$ws.Range("C3").Value = "camisa pessego"
$ws.Range("D3").Value = "'899.99"
$ws.Range("E3").Value = "'99.99"
$ws.Range("G3").Value = "Camisa-rosa"
$ws.Range("M3").Value = "'true"

# Restore the view/selection state left behind by the edit
$ws.Range("C3").Select() | Out-Null
